$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9: value changes from 107246 to 110245
$ws.Range("B9").Value = 110245

# C9: cell content removed entirely (was a validation status string)
$ws.Range("C9").ClearContents()

# Q9: value rounded from 456922.1424461872 to 456922
$ws.Range("Q9").Value = 456922

# R9: value rounded from 6200655.433266406 to 6200655
$ws.Range("R9").Value = 6200655

# Z9: cell content removed entirely (was "00:00")
$ws.Range("Z9").ClearContents()

# AB9: cell content removed entirely (was "00:00")
$ws.Range("AB9").ClearContents()
